$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 669.3913
$ws.Range("I15").Value = 669.3913
$ws.Range("K15").Value = 2008.1739
$ws.Range("M15").Value = -1839.1739

$ws.Range("H62").Value = 25002290
$ws.Range("I62").Value = 31251862
$ws.Range("K62").Value = 31251862
$ws.Range("M62").Value = -31251238

$ws.Range("H65").Value = 25002290
$ws.Range("I65").Value = 31251862
$ws.Range("K65").Value = 156259310
$ws.Range("M65").Value = -156256190

$ws.Range("H74").Value = 12179.357
$ws.Range("I74").Value = 11193.154
$ws.Range("K74").Value = 11193.154
$ws.Range("M74").Value = -10257.154

$ws.Range("H77").Value = 12179.357
$ws.Range("I77").Value = 11193.154
$ws.Range("K77").Value = 55965.77
$ws.Range("M77").Value = -51285.77

$ws.Range("H135").Value = 2409.1853
$ws.Range("I135").Value = 1135.5264
$ws.Range("K135").Value = 10219.7376
$ws.Range("M135").Value = -7684.7376

$ws.Range("H138").Value = 3934.7031
$ws.Range("I138").Value = 849.1818
$ws.Range("J138").Value = 5550.9287
$ws.Range("K138").Value = 2547.5454
$ws.Range("L138").Value = 16652.7861
$ws.Range("M138").Value = 2592.4546
$ws.Range("N138").Value = -26932.7861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 32667.576
$ws.Range("J2").Value = 4046.3333
$ws.Range("L2").Value = 4046.3333
$ws.Range("N2").Value = -4272.3333

$ws.Range("H32").Value = 3492.2964
$ws.Range("I32").Value = 3404.1177
$ws.Range("K32").Value = 3404.1177
$ws.Range("M32").Value = -3117.1177

$ws.Range("H57").Value = 9998
$ws.Range("I57").Value = 9998
$ws.Range("K57").Value = 9998
$ws.Range("M57").Value = -9514

$ws.Range("H74").Value = 24928.334
$ws.Range("J74").Value = 5098.7144
$ws.Range("L74").Value = 5098.7144
$ws.Range("N74").Value = -6846.7144

$ws.Range("H77").Value = 24928.334
$ws.Range("J77").Value = 5098.7144
$ws.Range("L77").Value = 25493.572
$ws.Range("N77").Value = -34229.572

$ws.Range("H116").Value = 32667.576
$ws.Range("J116").Value = 4046.3333
$ws.Range("L116").Value = 4046.3333
$ws.Range("N116").Value = -8634.3333

$ws.Range("H122").Value = 5160.316
$ws.Range("I122").Value = 4457.909
$ws.Range("J122").Value = 6126.125
$ws.Range("K122").Value = 13373.727
$ws.Range("L122").Value = 18378.375
$ws.Range("M122").Value = -10923.727
$ws.Range("N122").Value = -23278.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 32667.576
$ws.Range("J3").Value = 4046.3333
$ws.Range("L3").Value = 4046.3333
$ws.Range("N3").Value = -4274.3333

$ws.Range("H86").Value = 3076.75
$ws.Range("I86").Value = 1300
$ws.Range("K86").Value = 1300
$ws.Range("M86").Value = -177

$ws.Range("H89").Value = 3076.75
$ws.Range("I89").Value = 1300
$ws.Range("K89").Value = 6500
$ws.Range("M89").Value = -884

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3762.318
$ws.Range("J31").Value = 4070.2144
$ws.Range("L31").Value = 4070.2144
$ws.Range("N31").Value = -4660.2144

$ws.Range("H34").Value = 3762.318
$ws.Range("J34").Value = 4070.2144
$ws.Range("L34").Value = 4070.2144
$ws.Range("N34").Value = -4474.2144

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 175.25
$ws.Range("J15").Value = 200.33333
$ws.Range("L15").Value = 600.99999
$ws.Range("N15").Value = -880.99999

$ws.Range("H51").Value = 2213.125
$ws.Range("I51").Value = 1600
$ws.Range("J51").Value = 2417.5
$ws.Range("K51").Value = 4800
$ws.Range("L51").Value = 7252.5
$ws.Range("M51").Value = -4340
$ws.Range("N51").Value = -8172.5

$ws.Range("H57").Value = 11000
$ws.Range("I57").Value = 10000
$ws.Range("J57").Value = 12000
$ws.Range("K57").Value = 30000
$ws.Range("L57").Value = 36000
$ws.Range("M57").Value = -29441
$ws.Range("N57").Value = -37118

$ws.Range("H64").Value = 1000000000
$ws.Range("I64").Value = 1000000000
$ws.Range("K64").Value = 3000000000
$ws.Range("M64").Value = -2999999730

$ws.Range("H67").Value = 1000000000
$ws.Range("I67").Value = 1000000000
$ws.Range("K67").Value = 3000000000
$ws.Range("M67").Value = -2999999064

$ws.Range("H70").Value = 15000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 15000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H140").Value = 6053
$ws.Range("I140").Value = 6053
$ws.Range("K140").Value = 18159
$ws.Range("M140").Value = -12979

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 15874.417
$ws.Range("I55").Value = 14333.333
$ws.Range("J55").Value = 16388.111
$ws.Range("K55").Value = 14333.333
$ws.Range("L55").Value = 16388.111
$ws.Range("M55").Value = -14006.333
$ws.Range("N55").Value = -17042.111

$ws.Range("H135").Value = 99999.8
$ws.Range("J135").Value = 99999.8
$ws.Range("L135").Value = 99999.8
$ws.Range("N135").Value = -110139.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1528.85
$ws.Range("I22").Value = 812.4545000000001
$ws.Range("J22").Value = 2404.4443
$ws.Range("K22").Value = 812.4545000000001
$ws.Range("L22").Value = 2404.4443
$ws.Range("M22").Value = -517.4545000000001
$ws.Range("N22").Value = -2994.4443

$ws.Range("H27").Value = 1528.85
$ws.Range("I27").Value = 812.4545000000001
$ws.Range("J27").Value = 2404.4443
$ws.Range("K27").Value = 812.4545000000001
$ws.Range("L27").Value = 2404.4443
$ws.Range("M27").Value = -705.4545000000001
$ws.Range("N27").Value = -2618.4443

$ws.Range("H40").Value = 2504376
$ws.Range("I40").Value = 5001502
$ws.Range("K40").Value = 5001502
$ws.Range("M40").Value = -5001366

$ws.Range("H82").Value = 3075.2
$ws.Range("I82").Value = 3949.75
$ws.Range("J82").Value = 2492.1667
$ws.Range("K82").Value = 3949.75
$ws.Range("L82").Value = 2492.1667
$ws.Range("M82").Value = -3588.75
$ws.Range("N82").Value = -3214.1667

$ws.Range("H85").Value = 3075.2
$ws.Range("I85").Value = 3949.75
$ws.Range("J85").Value = 2492.1667
$ws.Range("K85").Value = 3949.75
$ws.Range("L85").Value = 2492.1667
$ws.Range("M85").Value = -2701.75
$ws.Range("N85").Value = -4988.1667

$ws.Range("H122").Value = 2105600.5
$ws.Range("I122").Value = 10000004
$ws.Range("J122").Value = 1228444.6
$ws.Range("K122").Value = 30000012
$ws.Range("L122").Value = 3685333.8
$ws.Range("M122").Value = -29997562
$ws.Range("N122").Value = -3690233.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6224.5
$ws.Range("I81").Value = 2499.5
$ws.Range("K81").Value = 4999
$ws.Range("M81").Value = -3938

$ws.Range("H84").Value = 6224.5
$ws.Range("I84").Value = 2499.5
$ws.Range("K84").Value = 24995
$ws.Range("M84").Value = -19691

$ws.Range("H136").Value = 1432251
$ws.Range("I136").Value = 1669625.4
$ws.Range("K136").Value = 5008876.199999999
$ws.Range("M136").Value = -5006326.199999999
